$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 'maternidade'
$ws.Range("B6").Value = 0.009843865097098409
$ws.Range("C6").Value = 0.1644615277481884
$ws.Range("D6").Value = $null
$ws.Range("E6").Value = $null
$ws.Range("A7").Value = 'linguagem e representação'
$ws.Range("B7").Value = 0.1149780503557664
$ws.Range("C7").Value = 0.06699772637601541
$ws.Range("D7").Value = $null
$ws.Range("E7").Value = $null
$ws.Range("A8").Value = 'família'
$ws.Range("B8").Value = 0.01180764123068177
$ws.Range("C8").Value = 0.217515846040044
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = $null
$ws.Range("A9").Value = 'resistência e luta'
$ws.Range("B9").Value = 0.09765825647091046
$ws.Range("C9").Value = 0.2060222982880944
$ws.Range("D9").Value = 0.6745230078563412
$ws.Range("E9").Value = 0
$ws.Range("A10").Value = 'saudade, luto ou perda'
$ws.Range("B10").Value = 0.7296010100853808
$ws.Range("C10").Value = 0.240792641016025
$ws.Range("D10").Value = 0.03010033444816049
$ws.Range("E10").Value = 0.006756756756756931
$ws.Range("A11").Value = 'sonho e fantasia'
$ws.Range("B11").Value = 0.6024657515707322
$ws.Range("C11").Value = 0.3660610585779649
$ws.Range("D11").Value = 0.03010033444816049
$ws.Range("E11").Value = 0.01342281879194638
$ws.Range("A12").Value = 'questão agrária e territorial'
$ws.Range("B12").Value = 0.007878426155611526
$ws.Range("C12").Value = 0.5534044896043896
$ws.Range("D12").Value = 0.7517123287671232
$ws.Range("E12").Value = 0.01342281879194638
$ws.Range("A13").Value = 'mobilidade'
$ws.Range("B13").Value = 0.1955561926033942
$ws.Range("C13").Value = 0.1733189057965303
$ws.Range("D13").Value = 0.575816674792784
$ws.Range("E13").Value = 0.01342281879194638
$ws.Range("A14").Value = 'tecnologia, inovação e sociedade'
$ws.Range("B14").Value = 0.009843865097098409
$ws.Range("C14").Value = 0.5554951447037371
$ws.Range("D14").Value = $null
$ws.Range("E14").Value = 0.01672240802675598
$ws.Range("A15").Value = 'sistema penitenciário'
$ws.Range("B15").Value = 0.01573021350765386
$ws.Range("C15").Value = 0.9261735546692708
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0.01672240802675598
$ws.Range("A16").Value = 'vida rural, vida no campo'
$ws.Range("B16").Value = 0.5404000685018602
$ws.Range("C16").Value = 0.393572857609351
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0.01672240802675598
$ws.Range("A17").Value = 'recreação, lazer e entretenimento'
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 0.07349522027797104
$ws.Range("D17").Value = 0.03010033444816049
$ws.Range("E17").Value = 0.01672240802675598
$ws.Range("A18").Value = 'estruturas sociais e econômicas'
$ws.Range("B18").Value = 0.3042982738151534
$ws.Range("C18").Value = 0.2851806637945615
$ws.Range("D18").Value = 0.02356902356902355
$ws.Range("E18").Value = 0.01672240802675598
$ws.Range("A19").Value = 'reflexão'
$ws.Range("B19").Value = 0.00197211133738422
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0.03010033444816049
$ws.Range("E19").Value = 0.1676627870399499
$ws.Range("A20").Value = 'vida cotidiana'
$ws.Range("B20").Value = 0.3622134994420639
$ws.Range("C20").Value = 0.2830586390250168
$ws.Range("D20").Value = 0.02684563758389252
$ws.Range("E20").Value = 0.2142038946162659
$ws.Range("A21").Value = 'cultural'
$ws.Range("B21").Value = 0.3908238931400018
$ws.Range("C21").Value = 0.316651444718541
$ws.Range("D21").Value = 0.6113642455105869
$ws.Range("E21").Value = 0.4020338983050848
$ws.Range("A22").Value = 'pandemia'
$ws.Range("B22").Value = 0.6374966323176273
$ws.Range("C22").Value = 0.5846045578451793
$ws.Range("D22").Value = 0.7447183098591549
$ws.Range("E22").Value = 0.438717067583047
$ws.Range("A23").Value = 'arte'
$ws.Range("B23").Value = 0.4473889112503763
$ws.Range("C23").Value = 0.1229787754676425
$ws.Range("D23").Value = 0.3614344133375275
$ws.Range("E23").Value = 0.4923873803170618
$ws.Range("A24").Value = 'memória e patrimônio'
$ws.Range("B24").Value = 0.6404480535751201
$ws.Range("C24").Value = 0.3870763395910517
$ws.Range("D24").Value = 0.2248949980908744
$ws.Range("E24").Value = 0.4965753424657535
$ws.Range("A25").Value = 'solidão'
$ws.Range("B25").Value = 0.4448519716756111
$ws.Range("C25").Value = 0.1615659732965542
$ws.Range("D25").Value = $null
$ws.Range("E25").Value = 0.4965753424657535
$ws.Range("A26").Value = 'dinâmica urbana'
$ws.Range("B26").Value = 0.5642947526978185
$ws.Range("C26").Value = 0.2428731689280502
$ws.Range("D26").Value = 0.3287037037037037
$ws.Range("E26").Value = 0.4965753424657535
$ws.Range("A27").Value = 'sonoridade e paisagem sonora'
$ws.Range("B27").Value = 0.3780716616176815
$ws.Range("C27").Value = 0.4662868848353157
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0.5033783783783785
$ws.Range("A28").Value = 'territorialidade e colonialismo'
$ws.Range("B28").Value = 0.7992186797845382
$ws.Range("C28").Value = 0.180562355026172
$ws.Range("D28").Value = 0.8892707140129821
$ws.Range("E28").Value = 0.5033783783783785
$ws.Range("A29").Value = 'política'
$ws.Range("B29").Value = 0.2658406082793572
$ws.Range("C29").Value = 0.8612121313295255
$ws.Range("D29").Value = 0.5101351351351351
$ws.Range("E29").Value = 0.5033783783783785
$ws.Range("A30").Value = 'saúde mental'
$ws.Range("B30").Value = 0.8391755387841752
$ws.Range("C30").Value = 0.3083490529686563
$ws.Range("D30").Value = 0.6206671026814912
$ws.Range("E30").Value = 0.5594405594405595
$ws.Range("A31").Value = 'poesia e ensaio'
$ws.Range("B31").Value = 0.3844201043045476
$ws.Range("C31").Value = 0.2480209576589355
$ws.Range("D31").Value = 0.4098494098494098
$ws.Range("E31").Value = 0.5800000000000002
$ws.Range("A32").Value = 'povos originários e comunidades tradicionais'
$ws.Range("B32").Value = 0.358805023395918
$ws.Range("C32").Value = 0.753079967989343
$ws.Range("D32").Value = 0.8586055582642613
$ws.Range("E32").Value = 0.6060037523452158
$ws.Range("A33").Value = 'mulher, feminino e feminismo'
$ws.Range("B33").Value = 0.741300314492549
$ws.Range("C33").Value = 0.4826546599346259
$ws.Range("D33").Value = 0.5438817238125196
$ws.Range("E33").Value = 0.656140350877193
$ws.Range("A34").Value = 'corpo, performance e expressão'
$ws.Range("B34").Value = 0.2398020297810572
$ws.Range("C34").Value = 0.65624667212239
$ws.Range("D34").Value = $null
$ws.Range("E34").Value = 0.6632302405498283
$ws.Range("A35").Value = 'distopia, ficção científica e futuros imaginados'
$ws.Range("B35").Value = 0.6236784203091036
$ws.Range("C35").Value = 0.3340710298986552
$ws.Range("D35").Value = 0.6745230078563412
$ws.Range("E35").Value = 0.6666666666666667
$ws.Range("A36").Value = 'violências e preconceitos de gênero'
$ws.Range("B36").Value = 0.3432649986348301
$ws.Range("C36").Value = 0.9715323261943718
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 0.6700336700336702
$ws.Range("A37").Value = 'sociedade e consumo'
$ws.Range("B37").Value = 0.2565198049716462
$ws.Range("C37").Value = 0.3817172538656722
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 0.6700336700336702
$ws.Range("A38").Value = 'moradia e habitação'
$ws.Range("B38").Value = 0.4501957803206739
$ws.Range("C38").Value = 0.1465016357325821
$ws.Range("D38").Value = 0.3287037037037037
$ws.Range("E38").Value = 0.6700336700336702
$ws.Range("A39").Value = 'violência'
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 0.09755202843081125
$ws.Range("D39").Value = 0.8033898305084746
$ws.Range("E39").Value = 0.6700336700336702
$ws.Range("A40").Value = 'trabalho e ofício'
$ws.Range("B40").Value = 0.6817643072219247
$ws.Range("C40").Value = 0.4818480702578282
$ws.Range("D40").Value = 0.6678121420389461
$ws.Range("E40").Value = 0.7482876712328766
$ws.Range("A41").Value = 'vida afetiva'
$ws.Range("B41").Value = 0.4046111906318235
$ws.Range("C41").Value = 0.2581704093286067
$ws.Range("D41").Value = 0.3614344133375275
$ws.Range("E41").Value = 0.7482876712328766
$ws.Range("A42").Value = 'alimentação e tratamentos tradicionais'
$ws.Range("B42").Value = 0.6407847426528892
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 0.4101694915254238
$ws.Range("E42").Value = 0.7482876712328766
$ws.Range("A43").Value = 'educação e socialização'
$ws.Range("B43").Value = 0.5344141900327821
$ws.Range("C43").Value = 0.1743901668865061
$ws.Range("D43").Value = 0.6
$ws.Range("E43").Value = 0.7482876712328766
$ws.Range("A44").Value = 'desinformação, populismo e polarização'
$ws.Range("B44").Value = 0.8487816303702602
$ws.Range("C44").Value = 0.9544999994540984
$ws.Range("D44").Value = 0.02684563758389252
$ws.Range("E44").Value = 0.8006779661016949
$ws.Range("A45").Value = 'ambiental'
$ws.Range("B45").Value = 0.8076672965220715
$ws.Range("C45").Value = 0.4572079115986321
$ws.Range("D45").Value = 0.5940959409594095
$ws.Range("E45").Value = 0.8077174623937214
$ws.Range("A46").Value = 'religião, espiritualidade e cosmologias'
$ws.Range("B46").Value = 0.8210032816280555
$ws.Range("C46").Value = 0.2800113589597854
$ws.Range("D46").Value = 0.8033898305084746
$ws.Range("E46").Value = 0.8877434135166095
$ws.Range("A47").Value = 'crises e desastres ambientais e sociais'
$ws.Range("B47").Value = 0.5350967631934763
$ws.Range("C47").Value = 0.3638649787863257
$ws.Range("D47").Value = 0.2551369863013699
$ws.Range("E47").Value = 0.9075180874488833
$ws.Range("A48").Value = 'pessoas com deficiência'
$ws.Range("B48").Value = 0.8339904252323289
$ws.Range("C48").Value = 0.4293576085057091
$ws.Range("D48").Value = 1
$ws.Range("E48").Value = 1
$ws.Range("A49").Value = 'gênero e sexualidade'
$ws.Range("B49").Value = 0.5116759821182142
$ws.Range("C49").Value = 0.9888511034491386
$ws.Range("D49").Value = 1
$ws.Range("E49").Value = 1
$ws.Range("A50").Value = 'direitos humanos'
$ws.Range("B50").Value = 0.6578256343361342
$ws.Range("C50").Value = 0.1149410793436468
$ws.Range("D50").Value = $null
$ws.Range("E50").Value = 1
$ws.Range("A51").Value = 'biografia'
$ws.Range("B51").Value = 0.2658406082793572
$ws.Range("C51").Value = 0.02624071623279086
$ws.Range("D51").Value = $null
$ws.Range("E51").Value = 1
$ws.Range("A52").Value = 'amizade'
$ws.Range("B52").Value = 0.6747223560389795
$ws.Range("C52").Value = 0.9725001544084381
$ws.Range("D52").Value = 0.03010033444816049
$ws.Range("E52").Value = 1
